$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, Dia, Mes, Ano, HoraInicio, HoraFim, QtdeHoras, LimiteHorasConsecRestante, HorasRestantes
$data = @(
  @(16, 1, 6, 2022, 19, 0, 5, 19, 91),
  @(17, 2, 6, 2022, 0, 7, 7, 12, 84),
  @(18, 2, 6, 2022, 19, 0, 5, 7, 79),
  @(19, 3, 6, 2022, 0, 7, 7, 0, 72),
  @(20, 6, 6, 2022, 0, 7, 7, 17, 65),
  @(21, 6, 6, 2022, 19, 0, 5, 12, 60),
  @(22, 7, 6, 2022, 0, 7, 7, 5, 53),
  @(23, 7, 6, 2022, 19, 0, 5, 0, 48),
  @(24, 9, 6, 2022, 0, 7, 7, 17, 41),
  @(25, 9, 6, 2022, 19, 0, 5, 12, 36),
  @(26, 10, 6, 2022, 0, 7, 7, 5, 29),
  @(27, 10, 6, 2022, 19, 0, 5, 0, 24),
  @(28, 13, 6, 2022, 0, 7, 7, 17, 17),
  @(29, 13, 6, 2022, 19, 0, 5, 12, 12),
  @(30, 14, 6, 2022, 0, 7, 7, 5, 5),
  @(31, 14, 6, 2022, 19, 0, 5, 0, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value2 = "Alicia"
    $ws.Cells.Item($r, 2).Value2 = "359.969.368-44"
    $ws.Cells.Item($r, 3).Value2 = "359.969.368-44"
    $ws.Cells.Item($r, 4).Value2 = $row[1]
    $ws.Cells.Item($r, 5).Value2 = $row[2]
    $ws.Cells.Item($r, 6).Value2 = $row[3]
    $ws.Cells.Item($r, 7).Value2 = $row[4]
    $ws.Cells.Item($r, 8).Value2 = $row[5]
    $ws.Cells.Item($r, 9).Value2 = $row[6]
    $ws.Cells.Item($r, 10).Value2 = $row[7]
    $ws.Cells.Item($r, 11).Value2 = $row[8]
}
